$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Rename the "Analysis 1" percent-change table rows:
#      pc_world          -> pc_rel_world
#      pc_world_v        -> pc_rel_world_v        (keep " (version for Power BI)")
#      pc_continent       -> pc_rel_continent       (drop " (version for Power BI)")
#      pc_continent_v     -> pc_rel_continent_v     (drop " (version for Power BI)")
#      pc_unregion        -> pc_rel_region          (drop " (version for Power BI)")
#      pc_unregion_v      -> pc_rel_region_v        (drop " (version for Power BI)")
# ---------------------------------------------------------------------

$p30 = $d.Paragraphs(30)
$p30.Range.Find.Execute("pc_world", $false, $false, $false, $false, $false, $true, 1, $false, "pc_rel_world", 2) | Out-Null

$p31 = $d.Paragraphs(31)
$p31.Range.Find.Execute("pc_world_v", $false, $false, $false, $false, $false, $true, 1, $false, "pc_rel_world_v", 2) | Out-Null

$p32 = $d.Paragraphs(32)
$p32.Range.Find.Execute("pc_continent", $false, $false, $false, $false, $false, $true, 1, $false, "pc_rel_continent", 2) | Out-Null
$p32b = $d.Paragraphs(32)
$p32b.Range.Find.Execute(" (version for Power BI)", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

$p33 = $d.Paragraphs(33)
$p33.Range.Find.Execute("pc_continent_v", $false, $false, $false, $false, $false, $true, 1, $false, "pc_rel_continent_v", 2) | Out-Null
$p33b = $d.Paragraphs(33)
$p33b.Range.Find.Execute(" (version for Power BI)", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

$p34 = $d.Paragraphs(34)
$p34.Range.Find.Execute("pc_unregion", $false, $false, $false, $false, $false, $true, 1, $false, "pc_rel_region", 2) | Out-Null
$p34b = $d.Paragraphs(34)
$p34b.Range.Find.Execute(" (version for Power BI)", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

$p35 = $d.Paragraphs(35)
$p35.Range.Find.Execute("pc_unregion_v", $false, $false, $false, $false, $false, $true, 1, $false, "pc_rel_region_v", 2) | Out-Null
$p35b = $d.Paragraphs(35)
$p35b.Range.Find.Execute(" (version for Power BI)", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. After "pc_rel_region_v" add a lone tab paragraph, then (keeping the
#    existing blank paragraph) three new UN-indicator style rows:
#      pc_ind_world
#      pc_ind_continent
#      pc_ind_region
# ---------------------------------------------------------------------

$p35 = $d.Paragraphs(35)
$p35.Range.InsertParagraphAfter()
$tabPara = $d.Paragraphs(36)
$tabPara.Range.InsertBefore("`t")

$blankPara = $d.Paragraphs(37)
$blankPara.Range.InsertParagraphAfter()
$indWorld = $d.Paragraphs(38)
$indWorld.Range.InsertBefore("`tpc_ind_world")

$indWorld = $d.Paragraphs(38)
$indWorld.Range.InsertParagraphAfter()
$indContinent = $d.Paragraphs(39)
$indContinent.Range.InsertBefore("`tpc_ind_continent")

$indContinent = $d.Paragraphs(39)
$indContinent.Range.InsertParagraphAfter()
$indRegion = $d.Paragraphs(40)
$indRegion.Range.InsertBefore("`tpc_ind_region")

# ---------------------------------------------------------------------
# 3. Add two extra blank paragraphs just before "Analysis 1 (...)"
# ---------------------------------------------------------------------

for ($i = 0; $i -lt 2; $i++) {
    $anchor = $d.Paragraphs(28)
    $anchor.Range.InsertParagraphBefore()
}

Write-Output "done"
